$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.570.31'
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").Value = '2.890.66'
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("E4").Value = '  +0.01%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'567.29"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  -4.66%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'142.88"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -4.00%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("D9").Value = '2.890.69'
$ws.Range("E9").Value = '  -2.05%  '
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'6.98"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = '  -2.45%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.145"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -3.51%  '
$ws.Range("E12").Value = '  -2.59%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  -2.05%  '
$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'31.83"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  -3.18%  '
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '3.373.30'
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("D17").Value = '61.617.05'
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("E18").Value = '  -2.39%  '
$ws.Range("D19").Value = '2.885.08'
$ws.Range("E19").Value = '  -2.47%  '
$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'432.03"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -2.26%  '
$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'13.04"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  -3.44%  '
$ws.Range("E22").Value = '  -1.53%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'79.39"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -1.93%  '
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'11.92"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("E26").Value = '  +0.02%  '
$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'9.90"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  -11.79%  '
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'2.00"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -6.07%  '
$ws.Range("E29").Value = '  +4.14%  '
$ws.Range("E30").Value = '  -4.59%  '
$ws.Range("E31").Value = '  -4.27%  '
$ws.Range("E32").Value = '  -8.78%  '
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("E34").Value = '  -2.40%  '
$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'25.51"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = '  -3.61%  '
$ws.Range("E36").Value = '  -3.11%  '
$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'48.87"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  -1.81%  '
$ws.Range("E39").Value = '  -5.87%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'2.80"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  -9.27%  '
$ws.Range("E41").Value = '  -3.54%  '
$ws.Range("E42").Value = '  -3.22%  '
$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'39.54"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("E44").Value = '  -5.15%  '
$ws.Range("D45").Value = '2.687.47'
$ws.Range("E45").Value = '  -0.79%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'132.90"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -2.01%  '
$ws.Range("E47").Value = '  -1.09%  '
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'346.07"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -4.07%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("E50").Value = '  -1.93%  '
$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'21.50"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -5.65%  '
